$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-16 already exist; the edit only touches columns A-I (J/K/L/M are
#     unchanged per the diff, so they are intentionally left untouched here). ---
$existingRows = @(
  @{ Row=2; A='DOmh6EPMeCY'; B='"What Would You Do?" — Bill O''Reilly Dissects Kamala Harris'' Reaction to Venezuela'; C='Subscribe to never miss an episode of No Spin News with Bill O''Reilly: https://www.youtube.com/channel/UC4OvD2yIbofl9l4dIlqSNMw'; D='2026-01-05T22:40:19Z'; E=135; F=31521; G=163; H='Bill O''Reilly'; I='https://www.youtube.com/watch?v=DOmh6EPMeCY' },
  @{ Row=3; A='5OJrS2xAoBQ'; B='Bill O''Reilly & Tim Graham on Internet Bias'; C='Newsbusters.org and Media Research Center''s Tim Graham joins Bill to discuss liberal bias on the internet.'; D='2026-01-04T13:00:18Z'; E=600; F=18106; G=74; H='Bill O''Reilly'; I='https://www.youtube.com/watch?v=5OJrS2xAoBQ' },
  @{ Row=4; A='s6W05lQNFuk'; B='He just accidentally EXPOSED Trump’s incompetence'; C='-- Marco Rubio tells Kristen Welker the United States is running Venezuela but cannot explain elections, governance, or the legal rationale in real time'; D='2026-01-06T03:30:48Z'; E=417; F=22171; G=277; H='David Pakman Show'; I='https://www.youtube.com/watch?v=s6W05lQNFuk' },
  @{ Row=5; A='UUrO5jWKlGU'; B='Trump is about to get REJECTED'; C='-- Donald Trump bets the Maduro seizure will unlock Venezuela’s oil, but allies, businesses, and markets resist a plan with no clear political transition'; D='2026-01-06T00:45:04Z'; E=343; F=155550; G=1166; H='David Pakman Show'; I='https://www.youtube.com/watch?v=UUrO5jWKlGU' },
  @{ Row=6; A='NJri0m1Uwiw'; B='OH NO: Trump has NO IDEA what’s going on'; C='-- Donald Trump admits he did not brief Congress yet says he spoke with oil companies and calls the Maduro seizure a kidnapping'; D='2026-01-05T23:30:16Z'; E=480; F=57721; G=593; H='David Pakman Show'; I='https://www.youtube.com/watch?v=NJri0m1Uwiw' },
  @{ Row=7; A='w0fPlxAgKaQ'; B='BOMBSHELL report: Trump health DETERIORATING QUICKLY'; C='-- The Wall Street Journal reports Donald Trump shows visible signs of decline, disputes doctors, and mislabels a CT scan as an MRI while aides manage optics'; D='2026-01-05T22:15:05Z'; E=522; F=232344; G=1882; H='David Pakman Show'; I='https://www.youtube.com/watch?v=w0fPlxAgKaQ' },
  @{ Row=8; A='qYk-K8lhZHc'; B='Trump invaded Venezuela because Maduro was dancing'; C='-- Reporting says Donald Trump escalates to military action after Nicolás Maduro dances on Venezuelan state television and aides treat it as personal mockery'; D='2026-01-05T21:01:04Z'; E=341; F=44235; G=745; H='David Pakman Show'; I='https://www.youtube.com/watch?v=qYk-K8lhZHc' },
  @{ Row=9; A='ldXNWL_w-Ao'; B='Is Trump even aware of what he just did?'; C='-- Donald Trump gives rambling, conflicting answers about Venezuela, oil, and military escalation that raise questions about his decision making'; D='2026-01-05T19:45:03Z'; E=487; F=307391; G=3225; H='David Pakman Show'; I='https://www.youtube.com/watch?v=ldXNWL_w-Ao' },
  @{ Row=10; A='sy8ncMFEPdE'; B='Trump invaded Venezuela… because Maduro danced? #shorts'; C='Become a Member: https://www.davidpakman.com/membership'; D='2026-01-05T17:21:46Z'; E=148; F=152418; G=1051; H='David Pakman Show'; I='https://www.youtube.com/watch?v=sy8ncMFEPdE' },
  @{ Row=11; A='1BU_O3mQkKI'; B='Trump gets BRUTALLY EXPOSED amid Venezuela invasion | Another Day'; C='Trump gets unwelcome SURPRISE amid Venezuela invasion | Another Day '; D='2026-01-06T05:00:02Z'; E=511; F=45972; G=454; H='Brian Tyler Cohen'; I='https://www.youtube.com/watch?v=1BU_O3mQkKI' },
  @{ Row=12; A='r81CCwvN5K8'; B='Trump makes GRAVE MISTAKE with Venezuela invasion'; C='INTERVIEW: Biden’s deputy national security adviser Jon Finer on Trump’s Venezuela invasion'; D='2026-01-06T02:59:07Z'; E=692; F=170989; G=1471; H='Brian Tyler Cohen'; I='https://www.youtube.com/watch?v=r81CCwvN5K8' },
  @{ Row=13; A='26fV3ovDlF8'; B='Mark Kelly drops BOMB on Pete Hegseth for trying to downgrade his rank'; C='INTERVIEW: Sen. Mark Kelly reacts to Pete Hegseth’s attempt to downrank him'; D='2026-01-06T00:59:05Z'; E=621; F=257354; G=1585; H='Brian Tyler Cohen'; I='https://www.youtube.com/watch?v=26fV3ovDlF8' },
  @{ Row=14; A='BhoTv296dp8'; B='How Huawei is winning the race for global 5G telecom dominance:  US firms didn''t even show up'; C='Half the world''s population now relies on telecommunications powered by Huawei and other Chinese companies.'; D='2026-01-06T02:14:01Z'; E=351; F=19760; G=308; H='Inside China Business'; I='https://www.youtube.com/watch?v=BhoTv296dp8' },
  @{ Row=15; A='6vw4hKAzo0k'; B='Revolutionary generator transforms Chinese factories into power plants'; C='Chinese engineers deployed the world''s first commercially viable sCO2 power generators, at a steel mill in Guizhou.'; D='2026-01-05T14:51:18Z'; E=281; F=57164; G=738; H='Inside China Business'; I='https://www.youtube.com/watch?v=6vw4hKAzo0k' },
  @{ Row=16; A='Bq1PDD5SWS0'; B='China plus Russia plus Iran plus North Korea: builds 70% of the world''s warships'; C='China''s dominance in commercial shipbuilding is hugely advantageous to the Chinese Navy, which is now the largest in the world.'; D='2026-01-04T12:09:37Z'; E=494; F=59070; G=749; H='Inside China Business'; I='https://www.youtube.com/watch?v=Bq1PDD5SWS0' }
)

foreach ($d in $existingRows) {
    $r = $d.Row
    $ws.Cells.Item($r, 1).Value = $d.A
    $ws.Cells.Item($r, 2).Value = $d.B
    $ws.Cells.Item($r, 3).Value = $d.C
    $ws.Cells.Item($r, 4).Value = $d.D
    $ws.Cells.Item($r, 5).Value = $d.E
    $ws.Cells.Item($r, 6).Value = $d.F
    $ws.Cells.Item($r, 7).Value = $d.G
    $ws.Cells.Item($r, 8).Value = $d.H
    $ws.Cells.Item($r, 9).Value = $d.I
}

# --- Rows 17-28 are newly appended; write every column A-M. ---
$newRows = @(
  @{ Row=17; A='qzue-WRW5MY'; B='China and Russia in the Arctic have NATO and Europe worried'; C='NATO and European officials are deeply concerned, as China and Russia enjoy scientific, commercial, and military breakthroughs across the Arctic region.'; D='2026-01-03T13:33:29Z'; E=513; F=40613; G=554; H='Inside China Business'; I='https://www.youtube.com/watch?v=qzue-WRW5MY'; J='en'; K='简体中文'; L=0; M='' },
  @{ Row=18; A='azAGvW91wss'; B='Aluminum prices soar on Trump tariffs, global shortages, and China supply chain moves'; C='American buyers of aluminum are paying record spreads over global benchmarks, amid Trump''s 50% tariffs and worldwide shortages of industrial metals.'; D='2026-01-02T13:46:45Z'; E=409; F=59704; G=481; H='Inside China Business'; I='https://www.youtube.com/watch?v=azAGvW91wss'; J='en'; K='简体中文'; L=0; M='' },
  @{ Row=19; A='t0a6IDHRq1o'; B='Nvidia is in big trouble, as Huawei rolls out 5G and AI across the world'; C='Nvidia faces severe challenges, as China''s monopolies on gallium allow its telecom providers to build low-cost 5G telecom across the world.  '; D='2025-12-30T11:00:55Z'; E=498; F=107288; G=1219; H='Inside China Business'; I='https://www.youtube.com/watch?v=t0a6IDHRq1o'; J='en'; K='简体中文'; L=0; M='' },
  @{ Row=20; A='qo8QnxqF92Y'; B='AI is revolutionizing Chinese coal production, and blowing up labor models everywhere else'; C='Coal prices are in steep decline across the world, and that should translate to collapsing profitability for coal miners.'; D='2025-12-29T13:26:55Z'; E=360; F=63274; G=586; H='Inside China Business'; I='https://www.youtube.com/watch?v=qo8QnxqF92Y'; J='en'; K='简体中文'; L=0; M='' },
  @{ Row=21; A='xged_Pzo35Q'; B='Trump VISIBLY STUNS Lindsey Graham in WEIRD RANT'; C='Donald Trump went on a weird rant which visibly stunned Lindsey Graham and also threatened to hurt blue states and cities.'; D='2026-01-05T19:39:43Z'; E=637; F=29394; G=341; H='Pondering Politics'; I='https://www.youtube.com/watch?v=xged_Pzo35Q'; J='en'; K='简体中文'; L=0; M='' },
  @{ Row=22; A='-3nzKAcEvxk'; B='Trump SINKS TO NEW LOW in DISGUSTING STUNT'; C='Donald Trump publicly mocks the brutal assassination of a Democratic lawmaker, whose children publicly beg him to stop.'; D='2026-01-05T18:34:22Z'; E=626; F=65437; G=576; H='Pondering Politics'; I='https://www.youtube.com/watch?v=-3nzKAcEvxk'; J='en'; K='简体中文'; L=0; M='' },
  @{ Row=23; A='Xj8JGZVgYFQ'; B='Marco Rubio ENDS CAREER by HUMILIATING Trump on LIVE TV'; C='Marco Rubio once again accidentally humiliated Donald Trump on live TV.'; D='2026-01-04T23:07:00Z'; E=693; F=245636; G=1773; H='Pondering Politics'; I='https://www.youtube.com/watch?v=Xj8JGZVgYFQ'; J='en'; K='简体中文'; L=0; M='' },
  @{ Row=24; A='4stBTanieHM'; B='🚨 Trump THREATENS MORE WAR as Venezuela GIVES HIM THE FINGER'; C='A furious Donald Trump publicly threatened to escalate his regime change war in Venezuela after the interim leader humiliated him.'; D='2026-01-04T19:07:00Z'; E=619; F=53202; G=721; H='Pondering Politics'; I='https://www.youtube.com/watch?v=4stBTanieHM'; J='en'; K='简体中文'; L=0; M='' },
  @{ Row=25; A='XUCBoQLmwyQ'; B='🚨 BOMBSHELL: Rubio SCREWS Trump, ADMITS DEADLY VENEZUELA LEAK'; C='Marco Rubio confirms bombshell reporting that the Trump administration accidentally leaked sensitive war information about their Venezuela war to the press.'; D='2026-01-04T18:07:00Z'; E=569; F=151417; G=953; H='Pondering Politics'; I='https://www.youtube.com/watch?v=XUCBoQLmwyQ'; J='en'; K='简体中文'; L=0; M='' },
  @{ Row=26; A='fLgZFbGt3KA'; B='🚨 STUNNING: Venezuela HUMILIATES Trump on LIVE TV'; C='Venezuela''s leadership, including interim president Delcy Rodriguez, publicly humiliated Donald Trump and rejected his claims they would roll over for him.'; D='2026-01-03T21:09:12Z'; E=632; F=687161; G=7457; H='Pondering Politics'; I='https://www.youtube.com/watch?v=fLgZFbGt3KA'; J='en'; K='简体中文'; L=0; M='' },
  @{ Row=27; A='FAn7bWn1uUM'; B='Will Trump Steal Greenland'; C='In this video I talk about what happened on the 4th of January! The craziest is Trump wanting Greenland and honestly not something I expected!'; D='2026-01-06T06:25:35Z'; E=666; F=1452; G=88; H='Omar Agamy'; I='https://www.youtube.com/watch?v=FAn7bWn1uUM'; J='en'; K='简体中文'; L=0; M='' },
  @{ Row=28; A='xCssF6vNGwU'; B='What Trump Did To Venezuela'; C='In this video I am talking about what the US did in Venezuela and how crazy the situation was in there. Honestly this story is enough to fill a whole day worth of videos cause it has been crazy,'; D='2026-01-05T04:57:58Z'; E=483; F=9622; G=381; H='Omar Agamy'; I='https://www.youtube.com/watch?v=xCssF6vNGwU'; J='en'; K='简体中文'; L=0; M='' }
)

foreach ($d in $newRows) {
    $r = $d.Row
    $ws.Cells.Item($r, 1).Value = $d.A
    $ws.Cells.Item($r, 2).Value = $d.B
    $ws.Cells.Item($r, 3).Value = $d.C
    $ws.Cells.Item($r, 4).Value = $d.D
    $ws.Cells.Item($r, 5).Value = $d.E
    $ws.Cells.Item($r, 6).Value = $d.F
    $ws.Cells.Item($r, 7).Value = $d.G
    $ws.Cells.Item($r, 8).Value = $d.H
    $ws.Cells.Item($r, 9).Value = $d.I
    $ws.Cells.Item($r, 10).Value = $d.J
    $ws.Cells.Item($r, 11).Value = $d.K
    $ws.Cells.Item($r, 12).Value = $d.L
    $ws.Cells.Item($r, 13).Value = $d.M
}

